# Retraining the model for Ulmeni
# Updates the Consumption_Forecast data (columns A & B) with new forecast
# values and extends the data range from row 89 to row 93.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(5590,5540,5500,5460,5440,5430,5420,5410,5400,5410,5420,5450,5480,5540,5610,5700,5830,6000,6220,6430,6640,6910,7180,7450,7650,7810,7970,8060,8120,8160,8170,8160,8140,8080,8000,7910,7830,7750,7660,7590,7520,7460,7390,7340,7330,7310,7290,7300,7320,7370,7410,7450,7530,7630,7700,7770,7840,7920,8000,8080,8170,8250,8310,8350,8370,8350,8340,8330,8300,8270,8240,8200,8140,8070,8000,7920,7810,7680,7580,7440,7250,7100,6970,6830,6690,6570,6480,6370,6240,6180,6140,6070)

$bValues = @(46006,46006.01041666666,46006.02083333334,46006.03125,46006.04166666666,46006.05208333334,46006.0625,46006.07291666666,46006.08333333334,46006.11458333334,46006.125,46006.13541666666,46006.14583333334,46006.15625,46006.16666666666,46006.17708333334,46006.1875,46006.19791666666,46006.20833333334,46006.21875,46006.22916666666,46006.23958333334,46006.25,46006.26041666666,46006.27083333334,46006.28125,46006.29166666666,46006.30208333334,46006.3125,46006.32291666666,46006.33333333334,46006.34375,46006.35416666666,46006.36458333334,46006.375,46006.38541666666,46006.39583333334,46006.40625,46006.41666666666,46006.42708333334,46006.4375,46006.44791666666,46006.45833333334,46006.46875,46006.47916666666,46006.48958333334,46006.5,46006.52083333334,46006.53125,46006.54166666666,46006.55208333334,46006.5625,46006.57291666666,46006.58333333334,46006.59375,46006.60416666666,46006.61458333334,46006.625,46006.63541666666,46006.64583333334,46006.65625,46006.66666666666,46006.67708333334,46006.6875,46006.69791666666,46006.71875,46006.72916666666,46006.73958333334,46006.75,46006.76041666666,46006.77083333334,46006.78125,46006.79166666666,46006.80208333334,46006.8125,46006.82291666666,46006.83333333334,46006.84375,46006.85416666666,46006.86458333334,46006.875,46006.88541666666,46006.89583333334,46006.90625,46006.91666666666,46006.92708333334,46006.9375,46006.94791666666,46006.95833333334,46006.96875,46006.97916666666,46006.98958333334)

# Style reference cell for column B (existing numFmt for timestamps)
$bStyleRange = $ws.Range("B2")

for ($i = 0; $i -lt $aValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = $bValues[$i]
    if ($row -gt 89) {
        $cellB.NumberFormat = $bStyleRange.NumberFormat
    }
}
